# Update the "RegDetails" sheet (sheet2): remove row 3 (testdata56/admin96)
# and change B2 from "admin95" to "bread88"; also update A2 to "data45".
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("RegDetails")

# Update row 2 values
$ws.Range("A2").Value = "data45"
$ws.Range("B2").Value = "bread88"

# Clear out row 3 entirely (was testdata56 / admin96)
$ws.Range("A3:B3").ClearContents()

# Set column B width to match new layout (~10.43 chars stored width)
$ws.Columns.Item(2).ColumnWidth = 9.6

# Move active selection to B2 to match new dimension
$ws.Activate()
$ws.Range("B2").Select()
